$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'26.366.10"
$ws.Range("E2").Value = '  +1.31%  '

# Row 3
$ws.Range("D3").Value = "'1.685.50"
$ws.Range("E3").Value = '  +1.04%  '

# Row 4
$ws.Range("D4").Value = "'1.008"
$ws.Range("E4").Value = '  +0.34%  '

# Row 5
$ws.Range("D5").Value = "'218.47"
$ws.Range("E5").Value = '  +0.80%  '

# Row 6
$ws.Range("D6").Value = "'0.5558"
$ws.Range("E6").Value = '  +9.10%  '

# Row 7
$ws.Range("D7").Value = "'1.008"
$ws.Range("E7").Value = '  +0.30%  '

# Row 8
$ws.Range("D8").Value = "'0.2710"
$ws.Range("E8").Value = '  +2.13%  '

# Row 9
$ws.Range("D9").Value = "'0.06510"
$ws.Range("E9").Value = '  +1.90%  '

# Row 10
$ws.Range("D10").Value = "'22.12"
$ws.Range("E10").Value = '  +1.57%  '

# Row 11
$ws.Range("D11").Value = "'0.07564"
$ws.Range("E11").Value = '  +1.48%  '

# Row 12
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").Value = "'4.554"
$ws.Range("E12").Value = '  +0.97%  '

# Row 13
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = "'1.684.74"
$ws.Range("E13").Value = '  +0.82%  '

# Row 14
$ws.Range("D14").Value = "'0.5816"
$ws.Range("E14").Value = '  -0.13%  '

# Row 15
$ws.Range("D15").Value = "'0.000008468"
$ws.Range("E15").Value = '  -0.92%  '

# Row 16
$ws.Range("D16").Value = "'65.47"
$ws.Range("E16").Value = '  +1.79%  '

# Row 17
$ws.Range("D17").Value = "'26.429.34"

# Row 18
$ws.Range("D18").Value = "'4.948"
$ws.Range("E18").Value = '  +0.23%  '

# Row 19
$ws.Range("D19").Value = "'1.008"
$ws.Range("E19").Value = '  +0.30%  '

# Row 20
$ws.Range("D20").Value = "'10.93"
$ws.Range("E20").Value = '  +1.55%  '

# Row 21
$ws.Range("D21").Value = "'191.63"
$ws.Range("E21").Value = '  +0.15%  '

# Row 22
$ws.Range("D22").Value = "'6.243"
$ws.Range("E22").Value = '  +0.73%  '

# Row 23
$ws.Range("E23").Value = '  +0.22%  '

# Row 24
$ws.Range("D24").Value = "'148.59"
$ws.Range("E24").Value = '  +2.68%  '

# Row 25
$ws.Range("D25").Value = "'0.1330"
$ws.Range("E25").Value = '  +11.01%  '

# Row 26
$ws.Range("D26").Value = "'7.909"
$ws.Range("E26").Value = '  +3.94%  '

# Row 27
$ws.Range("E27").Value = '  +1.36%  '

# Row 28
$ws.Range("D28").Value = "'0.06330"
$ws.Range("E28").Value = '  -3.84%  '

# Row 29
$ws.Range("D29").Value = "'1.398"
$ws.Range("E29").Value = '  +4.36%  '

# Row 30
$ws.Range("D30").Value = "'1.326"
$ws.Range("E30").Value = '  +0.74%  '

# Row 31
$ws.Range("D31").Value = "'3.603"
$ws.Range("E31").Value = '  +1.76%  '

# Row 32
$ws.Range("D32").Value = "'3.584"
$ws.Range("E32").Value = '  +2.05%  '

# Row 33
$ws.Range("D33").Value = "'1.672"
$ws.Range("E33").Value = '  +1.11%  '

# Row 34
$ws.Range("D34").Value = "'1.041"
$ws.Range("E34").Value = '  +2.26%  '

# Row 35
$ws.Range("D35").Value = "'0.6223"
$ws.Range("E35").Value = '  +1.71%  '

# Row 36
$ws.Range("D36").Value = "'2.402"
$ws.Range("E36").Value = '  +1.39%  '

# Row 37
$ws.Range("D37").Value = "'2.721"
$ws.Range("E37").Value = '  +1.18%  '

# Row 38
$ws.Range("D38").Value = "'6.239"
$ws.Range("E38").Value = '  -0.89%  '

# Row 39
$ws.Range("B39").Value = 'Maker'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D39").Value = "'1.116.85"
$ws.Range("E39").Value = '  +2.32%  '

# Row 40
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").Value = "'0.01632"
$ws.Range("E40").Value = '  +2.05%  '

# Row 41
$ws.Range("D41").Value = "'0.8773"
$ws.Range("E41").Value = '  +0.47%  '

# Row 42
$ws.Range("E42").Value = '  +0.42%  '

# Row 43
$ws.Range("D43").Value = "'100.64"
$ws.Range("E43").Value = '  -0.47%  '

# Row 44
$ws.Range("D44").Value = "'1.836.43"
$ws.Range("E44").Value = '  +1.12%  '

# Row 45
$ws.Range("B45").Value = 'BabyDogeCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D45").Value = "'0.00000000109"
$ws.Range("E45").Value = '  -4.15%  '

# Row 46
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = "'57.43"
$ws.Range("E46").Value = '  +1.95%  '

# Row 47
$ws.Range("D47").Value = "'8.240"
$ws.Range("E47").Value = '  +2.27%  '

# Row 48
$ws.Range("D48").Value = "'1.005"
$ws.Range("E48").Value = '  -0.18%  '

# Row 49
$ws.Range("D49").Value = "'0.05279"
$ws.Range("E49").Value = '  +0.93%  '

# Row 50
$ws.Range("B50").Value = 'Mantle'
$ws.Range("C50").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D50").Value = "'0.4300"
$ws.Range("E50").Value = '  +0.25%  '

# Row 51
$ws.Range("B51").Value = 'Aptos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D51").Value = "'6.092"
$ws.Range("E51").Value = '  +1.12%  '
